$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as TEXT (matches the source data,
# which keeps price/volume figures as inline strings, not numbers),
# then restore the cell to its original unstyled state so no stray
# number-format / style footprint is left behind.
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.NumberFormat = "General"
    $Range.Style = "Normal"
}

# --- Update Price (D) and Volume(1h) (E) columns for rows 2-48 ---
$ws.Range("D2").Value = '26.183.23'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.590.43'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws.Range("D5") '211.47'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  -0.45%  '
Set-TextValue $ws.Range("D9") '0.0604'
$ws.Range("E9").Value = '  +0.17%  '
Set-TextValue $ws.Range("D10") '18.83'
$ws.Range("E10").Value = '  -1.86%  '
Set-TextValue $ws.Range("D11") '0.0853'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").Value = '1.821.54'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = '1.652.01'
$ws.Range("E13").Value = '  +4.32%  '
Set-TextValue $ws.Range("D14") '4.00'
$ws.Range("E14").Value = '  -0.21%  '
Set-TextValue $ws.Range("D15") '0.503'
$ws.Range("E15").Value = '  -2.52%  '
Set-TextValue $ws.Range("D16") '63.44'
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = '26.161.33'
$ws.Range("E17").Value = '  -0.12%  '
Set-TextValue $ws.Range("D18") '227.39'
$ws.Range("E18").Value = '  +6.21%  '
$ws.Range("D19").Value = '0.0₃0718'
$ws.Range("E19").Value = '  -0.77%  '
Set-TextValue $ws.Range("D20") '7.53'
$ws.Range("E20").Value = '  +3.64%  '
$ws.Range("E21").Value = '  +0.06%  '
Set-TextValue $ws.Range("D22") '4.22'
$ws.Range("E22").Value = '  -0.43%  '
Set-TextValue $ws.Range("D23") '2.16'
$ws.Range("E23").Value = '  +0.88%  '
Set-TextValue $ws.Range("D24") '8.87'
$ws.Range("E24").Value = '  -0.95%  '
Set-TextValue $ws.Range("D25") '145.32'
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("E26").Value = '  +0.11%  '
Set-TextValue $ws.Range("D27") '6.93'
$ws.Range("E27").Value = '  -0.87%  '
Set-TextValue $ws.Range("D28") '0.112'
$ws.Range("E28").Value = '  +0.32%  '
Set-TextValue $ws.Range("D29") '15.29'
$ws.Range("E29").Value = '  +1.37%  '
Set-TextValue $ws.Range("D30") '0.0491'
$ws.Range("E30").Value = '  -0.86%  '
Set-TextValue $ws.Range("D31") '1.15'
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("D33").Value = '1.448.84'
$ws.Range("E33").Value = '  +3.63%  '
$ws.Range("E34").Value = '  +0.17%  '
Set-TextValue $ws.Range("D35") '2.43'
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("E36").Value = '  +0.15%  '
Set-TextValue $ws.Range("D37") '0.562'
$ws.Range("E37").Value = '  -4.10%  '
$ws.Range("E38").Value = '  -1.30%  '
Set-TextValue $ws.Range("D39") '0.815'
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("E41").Value = '  +0.18%  '
Set-TextValue $ws.Range("D42") '2.17'
$ws.Range("E42").Value = '  +1.51%  '
Set-TextValue $ws.Range("D43") '0.926'
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").Value = '1.732.07'
$ws.Range("E44").Value = '  +0.60%  '
Set-TextValue $ws.Range("D45") '0.754'
$ws.Range("E45").Value = '  -1.46%  '
Set-TextValue $ws.Range("D46") '60.15'
$ws.Range("E46").Value = '  -1.30%  '
Set-TextValue $ws.Range("D47") '87.41'
$ws.Range("E47").Value = '  +1.95%  '
Set-TextValue $ws.Range("D48") '1.47'
$ws.Range("E48").Value = '  -0.71%  '

# --- Rows 49-51 restructuring ---
# A new coin (BabyDogeCoin) was inserted at row 49, pushing the former
# row 49 (Cronos) down to row 50 and the former row 50 (USDD) down to
# row 51; the former row 51 (EnergySwap) drops off the bottom of the list.
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D49") '0.0₇0976'
$ws.Range("E49").Value = '  -6.49%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '0.0500'
$ws.Range("E50").Value = '  +0.05%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue $ws.Range("D51") '1.00'
$ws.Range("E51").Value = '  +0.07%  '
